# Update the EPEX Spot prices workbook with the latest day of data.
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add a new day column (K) with header "24-jun" ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Header cell, matching style of the other date headers (bold/centered/bordered)
$ws1.Range("J1").Copy()
$ws1.Range("K1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("K1").Value = "24-jun"

# Hourly price values for 24-jun
$k1Values = @(
    81.34,
    66.17,
    63.79,
    45.14,
    46.62,
    71.3,
    84.53,
    90.09,
    73.55,
    37.04,
    1.34,
    0,
    -0.01,
    -0.01,
    -0.03,
    -0.01,
    7.5,
    20.06,
    71.95,
    104.61,
    125.4,
    125.49,
    131.91,
    102.48
)

$row = 2
foreach ($val in $k1Values) {
    $ws1.Cells.Item($row, 11).Value = $val
    $row++
}

# --- Sheet "Gaz": append the 2025-06-23 price ---
$ws2 = $wb.Worksheets.Item("Gaz")
# Force the date-shaped text to stay plain text (not get auto-converted to a
# date serial number), then drop the number format again so the cell ends up
# unstyled just like the existing date cells above it.
$ws2.Range("A7").NumberFormat = "@"
$ws2.Range("A7").Value = "2025-06-23"
$ws2.Range("A7").ClearFormats()
$ws2.Range("B7").Value = 40.9

# --- Sheet "CO2": append the 2025-06-23 price ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A7").NumberFormat = "@"
$ws3.Range("A7").Value = "2025-06-23"
$ws3.Range("A7").ClearFormats()
$ws3.Range("B7").Value = 71.88
